$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "65.989.59"
Set-TextValue "E2" "  +1.00%  "
Set-TextValue "D3" "3.313.38"
Set-TextValue "E3" "  +0.50%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.22%  "
Set-TextValue "D5" "188.31"
Set-TextValue "E5" "  +5.18%  "
Set-TextValue "D6" "556.84"
Set-TextValue "E6" "  +0.40%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.05%  "
Set-TextValue "E8" "  -0.30%  "
Set-TextValue "D9" "3.308.63"
Set-TextValue "E9" "  +0.75%  "
Set-TextValue "D10" "0.184"
Set-TextValue "E10" "  +0.73%  "
Set-TextValue "E11" "  +1.08%  "
Set-TextValue "D12" "47.47"
Set-TextValue "E12" "  +1.27%  "
Set-TextValue "E13" "  +3.79%  "
Set-TextValue "E14" "  +2.52%  "
Set-TextValue "D15" "3.839.52"
Set-TextValue "E15" "  +0.25%  "
Set-TextValue "D16" "607.60"
Set-TextValue "E16" "  +2.08%  "
Set-TextValue "B17" "Chainlink"
Set-TextValue "C17" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D17" "18.04"
Set-TextValue "E17" "  +0.57%  "
Set-TextValue "B18" "WrappedBTC"
Set-TextValue "C18" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D18" "66.015.91"
Set-TextValue "E18" "  +1.13%  "
Set-TextValue "E19" "  +1.24%  "
Set-TextValue "D20" "3.328.24"
Set-TextValue "E20" "  +1.03%  "
Set-TextValue "D21" "11.13"
Set-TextValue "E21" "  -1.74%  "
Set-TextValue "D22" "0.909"
Set-TextValue "E22" "  +1.59%  "
Set-TextValue "D23" "18.43"
Set-TextValue "E23" "  +9.06%  "
Set-TextValue "E24" "  +0.88%  "
Set-TextValue "D25" "100.19"
Set-TextValue "E25" "  -1.49%  "
Set-TextValue "D26" "3.97"
Set-TextValue "E26" "  +0.43%  "
Set-TextValue "E27" "  +4.79%  "
Set-TextValue "D28" "5.94"
Set-TextValue "E28" "  -0.64%  "
Set-TextValue "D29" "9.61"
Set-TextValue "E29" "  +4.51%  "
Set-TextValue "D30" "8.71"
Set-TextValue "E30" "  +1.32%  "
Set-TextValue "D31" "30.38"
Set-TextValue "E31" "  +0.00%  "
Set-TextValue "D32" "6.78"
Set-TextValue "E32" "  +9.66%  "
Set-TextValue "D33" "3.87"
Set-TextValue "E33" "  +1.63%  "
Set-TextValue "D34" "581.42"
Set-TextValue "E34" "  +13.07%  "
Set-TextValue "D35" "11.11"
Set-TextValue "E35" "  +1.44%  "
Set-TextValue "E36" "  +1.57%  "
Set-TextValue "E37" "  +0.17%  "
Set-TextValue "D38" "57.04"
Set-TextValue "E38" "  -0.05%  "
Set-TextValue "D39" "3.703.93"
Set-TextValue "E39" "  -1.82%  "
Set-TextValue "B40" "InjectiveProtocol"
Set-TextValue "C40" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D40" "34.09"
Set-TextValue "E40" "  +7.56%  "
Set-TextValue "B41" "PEPE"
Set-TextValue "C41" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D41" "0.0₃0722"
Set-TextValue "E41" "  +2.68%  "
Set-TextValue "B42" "Kaspa"
Set-TextValue "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.130"
Set-TextValue "E42" "  +5.59%  "
Set-TextValue "B43" "CoreDAO"
Set-TextValue "C43" "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D43" "3.44"
Set-TextValue "E43" "  +12.54%  "
Set-TextValue "D44" "3.29"
Set-TextValue "E44" "  -4.36%  "
Set-TextValue "E45" "  +2.31%  "
Set-TextValue "D46" "0.342"
Set-TextValue "E46" "  +1.76%  "
Set-TextValue "E47" "  +3.07%  "
Set-TextValue "D48" "0.0422"
Set-TextValue "E48" "  +3.54%  "
Set-TextValue "E49" "  +1.09%  "
Set-TextValue "E50" "  +0.29%  "
Set-TextValue "E51" "  -0.20%  "
